# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending list of "Periodo Mora" values (1703..1812, 1901..1912, 2001..2003)
# replacing the previous descending list (2003..1703).
$periods = @(
    "1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i

    # Column E: Periodo Mora (text, matches shared-string period code)
    $ws.Cells.Item($row, 5).Value = $periods[$i]

    # Column F: Valor Mora - 29509 for periods up to 1808 (row 33), 31249 from 1809 (row 34) onward
    if ($row -le 33) {
        $ws.Cells.Item($row, 6).Value = 29509
    } else {
        $ws.Cells.Item($row, 6).Value = 31249
    }

    # Column G: Salario Basico - updated for every row
    $ws.Cells.Item($row, 7).Value = 781242
}

# Move the logo image to the left (matches updated anchor offsets in drawing1.xml).
# Use a precise absolute position (19pt left of the original) rather than a relative
# "-19" offset, since the relative form compounds the Left-property's display rounding.
$shp = $ws.Shapes.Item(1)
$shp.Left = 56.41281680610236
